$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1711.8846
$ws.Range("I98").Value = 1204.9546
$ws.Range("J98").Value = 4500
$ws.Range("K98").Value = 1204.9546
$ws.Range("L98").Value = 4500
$ws.Range("M98").Value = 293.0454
$ws.Range("N98").Value = -7496
$ws.Range("H112").Value = 11365204
$ws.Range("J112").Value = 13159026
$ws.Range("L112").Value = 39477078
$ws.Range("N112").Value = -39479294
$ws.Range("H122").Value = 1711.8846
$ws.Range("I122").Value = 1204.9546
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 3614.8638
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -1164.8638
$ws.Range("N122").Value = -18400
$ws.Range("H131").Value = 2726
$ws.Range("I131").Value = 2780.3333
$ws.Range("J131").Value = 2400
$ws.Range("K131").Value = 8340.999899999999
$ws.Range("L131").Value = 7200
$ws.Range("M131").Value = -3300.999899999999
$ws.Range("N131").Value = -17280
$ws.Range("H137").Value = 3230118.2
$ws.Range("I137").Value = 4004071.2
$ws.Range("J137").Value = 5314.8335
$ws.Range("K137").Value = 12012213.6
$ws.Range("L137").Value = 15944.5005
$ws.Range("M137").Value = -12009663.6
$ws.Range("N137").Value = -21044.5005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3200
$ws.Range("I61").Value = 2141.6667
$ws.Range("J61").Value = 6375
$ws.Range("K61").Value = 2141.6667
$ws.Range("L61").Value = 6375
$ws.Range("M61").Value = -1929.6667
$ws.Range("N61").Value = -6799
$ws.Range("H74").Value = 1466.75
$ws.Range("I74").Value = 1551.6364
$ws.Range("J74").Value = 1280
$ws.Range("K74").Value = 1551.6364
$ws.Range("L74").Value = 1280
$ws.Range("M74").Value = -677.6364000000001
$ws.Range("N74").Value = -3028
$ws.Range("H77").Value = 1466.75
$ws.Range("I77").Value = 1551.6364
$ws.Range("J77").Value = 1280
$ws.Range("K77").Value = 7758.182000000001
$ws.Range("L77").Value = 6400
$ws.Range("M77").Value = -3390.182000000001
$ws.Range("N77").Value = -15136
$ws.Range("H136").Value = 3200
$ws.Range("I136").Value = 2141.6667
$ws.Range("J136").Value = 6375
$ws.Range("K136").Value = 6425.000100000001
$ws.Range("L136").Value = 19125
$ws.Range("M136").Value = -3875.000100000001
$ws.Range("N136").Value = -24225

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 292.75
$ws.Range("I22").Value = 292.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 292.75
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -119.75
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 537.5294
$ws.Range("I94").Value = 531.2857
$ws.Range("J94").Value = 566.6667
$ws.Range("K94").Value = 531.2857
$ws.Range("L94").Value = 566.6667
$ws.Range("M94").Value = -80.28570000000002
$ws.Range("N94").Value = -1468.6667
$ws.Range("H99").Value = 3801
$ws.Range("I99").Value = 2681.6
$ws.Range("K99").Value = 2681.6
$ws.Range("M99").Value = -1183.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 57008
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 70010
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 70010
$ws.Range("M23").Value = -4760
$ws.Range("N23").Value = -70490
$ws.Range("H27").Value = 57008
$ws.Range("I27").Value = 5000
$ws.Range("J27").Value = 70010
$ws.Range("K27").Value = 5000
$ws.Range("L27").Value = 70010
$ws.Range("M27").Value = -4808
$ws.Range("N27").Value = -70394
$ws.Range("H99").Value = 2217.6875
$ws.Range("I99").Value = 1641.3
$ws.Range("J99").Value = 3178.3333
$ws.Range("K99").Value = 1641.3
$ws.Range("L99").Value = 3178.3333
$ws.Range("M99").Value = -143.3
$ws.Range("N99").Value = -6174.3333
$ws.Range("H126").Value = 2217.6875
$ws.Range("I126").Value = 1641.3
$ws.Range("J126").Value = 3178.3333
$ws.Range("K126").Value = 4923.9
$ws.Range("L126").Value = 9534.999899999999
$ws.Range("M126").Value = -2453.9
$ws.Range("N126").Value = -14474.9999
$ws.Range("H132").Value = 1974.9546
$ws.Range("I132").Value = 1850.2174
$ws.Range("J132").Value = 2111.5715
$ws.Range("K132").Value = 5550.6522
$ws.Range("L132").Value = 6334.7145
$ws.Range("M132").Value = -3020.6522
$ws.Range("N132").Value = -11394.7145

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2780
$ws.Range("J20").Value = 3250
$ws.Range("L20").Value = 9750
$ws.Range("N20").Value = -10204
$ws.Range("H123").Value = 3550
$ws.Range("I123").Value = 500
$ws.Range("K123").Value = 1500
$ws.Range("M123").Value = 950
$ws.Range("H129").Value = 26597.38
$ws.Range("I129").Value = 4103.75
$ws.Range("J129").Value = 40439.617
$ws.Range("K129").Value = 12311.25
$ws.Range("L129").Value = 121318.851
$ws.Range("M129").Value = -7311.25
$ws.Range("N129").Value = -131318.851
$ws.Range("H130").Value = 2399.8572
$ws.Range("J130").Value = 2399.8572
$ws.Range("L130").Value = 7199.571599999999
$ws.Range("N130").Value = -17239.5716
$ws.Range("H131").Value = 2166.182
$ws.Range("J131").Value = 2152.8462
$ws.Range("L131").Value = 6458.5386
$ws.Range("N131").Value = -16538.5386
$ws.Range("H133").Value = 5603.067
$ws.Range("I133").Value = 6581.6665
$ws.Range("J133").Value = 4950.6665
$ws.Range("K133").Value = 19744.9995
$ws.Range("L133").Value = 14851.9995
$ws.Range("M133").Value = -14684.9995
$ws.Range("N133").Value = -24971.9995
$ws.Range("H134").Value = 3186.7896
$ws.Range("I134").Value = 1575
$ws.Range("J134").Value = 4977.6665
$ws.Range("K134").Value = 4725
$ws.Range("L134").Value = 14932.9995
$ws.Range("M134").Value = 345
$ws.Range("N134").Value = -25072.9995
$ws.Range("H136").Value = 1779.1428
$ws.Range("J136").Value = 2496
$ws.Range("L136").Value = 7488
$ws.Range("N136").Value = -17688
$ws.Range("H137").Value = 3137.3044
$ws.Range("I137").Value = 2321.6667
$ws.Range("J137").Value = 4027.0908
$ws.Range("K137").Value = 6965.000100000001
$ws.Range("L137").Value = 12081.2724
$ws.Range("M137").Value = -1865.000100000001
$ws.Range("N137").Value = -22281.2724
$ws.Range("H138").Value = 4473.222
$ws.Range("J138").Value = 7399.8
$ws.Range("L138").Value = 22199.4
$ws.Range("N138").Value = -32479.4
$ws.Range("H139").Value = 4371.1177
$ws.Range("I139").Value = 3313.75
$ws.Range("J139").Value = 5311
$ws.Range("K139").Value = 9941.25
$ws.Range("L139").Value = 15933
$ws.Range("M139").Value = -4801.25
$ws.Range("N139").Value = -26213
$ws.Range("H140").Value = 2735.625
$ws.Range("I140").Value = 441.1111
$ws.Range("J140").Value = 5685.7144
$ws.Range("K140").Value = 1323.3333
$ws.Range("L140").Value = 17057.1432
$ws.Range("M140").Value = 3856.6667
$ws.Range("N140").Value = -27417.1432
$ws.Range("H141").Value = 3260
$ws.Range("I141").Value = 3433.3333
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 10299.9999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -5119.999899999999
$ws.Range("N141").Value = -19360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 30667.143
$ws.Range("I102").Value = 1765.2858
$ws.Range("K102").Value = 1765.2858
$ws.Range("M102").Value = -143.2858000000001
$ws.Range("H122").Value = 3986.6428
$ws.Range("I122").Value = 2654.4707
$ws.Range("K122").Value = 7963.4121
$ws.Range("M122").Value = -5513.4121
$ws.Range("H126").Value = 3174.5
$ws.Range("I126").Value = 1781.6666
$ws.Range("J126").Value = 3771.4285
$ws.Range("K126").Value = 5344.9998
$ws.Range("L126").Value = 11314.2855
$ws.Range("M126").Value = -2874.9998
$ws.Range("N126").Value = -16254.2855
$ws.Range("H137").Value = 63516.363
$ws.Range("J137").Value = 63516.363
$ws.Range("L137").Value = 63516.363
$ws.Range("N137").Value = -73716.363
$ws.Range("H138").Value = 30000
$ws.Range("J138").Value = 30000
$ws.Range("L138").Value = 30000
$ws.Range("N138").Value = -40280

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1878.25
$ws.Range("I7").Value = 1179.4286
$ws.Range("K7").Value = 1179.4286
$ws.Range("M7").Value = -1067.4286
$ws.Range("H40").Value = 8269.857
$ws.Range("I40").Value = 7981.5
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 7981.5
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -7845.5
$ws.Range("N40").Value = -10272
$ws.Range("H122").Value = 2828.762
$ws.Range("I122").Value = 2416.5833
$ws.Range("J122").Value = 3378.3333
$ws.Range("K122").Value = 7249.749899999999
$ws.Range("L122").Value = 10134.9999
$ws.Range("M122").Value = -4799.749899999999
$ws.Range("N122").Value = -15034.9999
$ws.Range("H126").Value = 1878.25
$ws.Range("I126").Value = 1179.4286
$ws.Range("K126").Value = 3538.2858
$ws.Range("M126").Value = -1068.2858
$ws.Range("H132").Value = 2287.093
$ws.Range("I132").Value = 1553.25
$ws.Range("J132").Value = 2721.963
$ws.Range("K132").Value = 4659.75
$ws.Range("L132").Value = 8165.889000000001
$ws.Range("M132").Value = -2129.75
$ws.Range("N132").Value = -13225.889

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 557850.6
$ws.Range("I122").Value = 771078.9
$ws.Range("K122").Value = 2313236.7
$ws.Range("M122").Value = -2310786.7
$ws.Range("H132").Value = 242667.6
$ws.Range("I132").Value = 401283.1
$ws.Range("K132").Value = 1203849.3
$ws.Range("M132").Value = -1201319.3
